$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("H58").Value = 2182.4092
$ws.Range("I58").Value = 230.61539
$ws.Range("J58").Value = 5001.6665
$ws.Range("K58").Value = 691.84617
$ws.Range("L58").Value = 15004.9995
$ws.Range("M58").Value = -541.84617
$ws.Range("N58").Value = -15304.9995
$ws.Range("H92").Value = 7283.2856
$ws.Range("I92").Value = 4830.5
$ws.Range("K92").Value = 4830.5
$ws.Range("M92").Value = -3582.5
$ws.Range("H96").Value = 322.64706
$ws.Range("I96").Value = 314.6154
$ws.Range("J96").Value = 348.75
$ws.Range("K96").Value = 943.8462000000001
$ws.Range("L96").Value = 1046.25
$ws.Range("M96").Value = 429.1537999999999
$ws.Range("N96").Value = -3792.25
$ws.Range("H100").Value = 2338.4
$ws.Range("I100").Value = 564.3333
$ws.Range("J100").Value = 4999.5
$ws.Range("K100").Value = 564.3333
$ws.Range("L100").Value = 4999.5
$ws.Range("M100").Value = -23.33330000000001
$ws.Range("N100").Value = -6081.5
$ws.Range("H101").Value = 66969.164
$ws.Range("I101").Value = 921.44446
$ws.Range("J101").Value = 265112.34
$ws.Range("K101").Value = 2764.33338
$ws.Range("L101").Value = 795337.02
$ws.Range("M101").Value = -1142.33338
$ws.Range("N101").Value = -798581.02
$ws.Range("H103").Value = 890.1111
$ws.Range("I103").Value = 1018.4
$ws.Range("J103").Value = 840.7692
$ws.Range("K103").Value = 3055.2
$ws.Range("L103").Value = 2522.3076
$ws.Range("M103").Value = -2469.2
$ws.Range("N103").Value = -3694.3076
$ws.Range("H106").Value = 3109.9
$ws.Range("I106").Value = 2639.8
$ws.Range("J106").Value = 3580
$ws.Range("K106").Value = 2639.8
$ws.Range("L106").Value = 3580
$ws.Range("M106").Value = -2008.8
$ws.Range("N106").Value = -4842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 95188.60000000001
$ws.Range("J23").Value = 94876.22
$ws.Range("L23").Value = 94876.22
$ws.Range("N23").Value = -95394.22
$ws.Range("H32").Value = 6429.9287
$ws.Range("I32").Value = 1564.1794
$ws.Range("J32").Value = 17592.53
$ws.Range("K32").Value = 1564.1794
$ws.Range("L32").Value = 17592.53
$ws.Range("M32").Value = -1277.1794
$ws.Range("N32").Value = -18166.53
$ws.Range("H33").Value = 17500
$ws.Range("J33").Value = 17500
$ws.Range("L33").Value = 17500
$ws.Range("N33").Value = -18158
$ws.Range("H61").Value = 48236.684
$ws.Range("I61").Value = 2872.9
$ws.Range("J61").Value = 501874.5
$ws.Range("K61").Value = 2872.9
$ws.Range("L61").Value = 501874.5
$ws.Range("M61").Value = -2660.9
$ws.Range("N61").Value = -502298.5
$ws.Range("H97").Value = 1114.2106
$ws.Range("I97").Value = 869.5
$ws.Range("J97").Value = 1799.4
$ws.Range("K97").Value = 869.5
$ws.Range("L97").Value = 1799.4
$ws.Range("M97").Value = -373.5
$ws.Range("N97").Value = -2791.4
$ws.Range("H132").Value = 2296.3872
$ws.Range("I132").Value = 2188.7307
$ws.Range("J132").Value = 2856.2
$ws.Range("K132").Value = 6566.1921
$ws.Range("L132").Value = 8568.599999999999
$ws.Range("M132").Value = -4036.1921
$ws.Range("N132").Value = -13628.6
$ws.Range("H136").Value = 48236.684
$ws.Range("I136").Value = 2872.9
$ws.Range("J136").Value = 501874.5
$ws.Range("K136").Value = 8618.700000000001
$ws.Range("L136").Value = 1505623.5
$ws.Range("M136").Value = -6068.700000000001
$ws.Range("N136").Value = -1510723.5
$ws.Range("H139").Value = 86283.8
$ws.Range("J139").Value = 86283.8
$ws.Range("L139").Value = 86283.8
$ws.Range("N139").Value = -96563.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 44856.5
$ws.Range("I20").Value = 70077
$ws.Range("J20").Value = 2822.3333
$ws.Range("K20").Value = 70077
$ws.Range("L20").Value = 2822.3333
$ws.Range("M20").Value = -69830
$ws.Range("N20").Value = -3316.3333
$ws.Range("H81").Value = 27390
$ws.Range("J81").Value = 27390
$ws.Range("L81").Value = 27390
$ws.Range("N81").Value = -29512
$ws.Range("H84").Value = 27390
$ws.Range("J84").Value = 27390
$ws.Range("L84").Value = 82170
$ws.Range("N84").Value = -92778
$ws.Range("H94").Value = 1790.091
$ws.Range("I94").Value = 1943
$ws.Range("J94").Value = 1606.6
$ws.Range("K94").Value = 1943
$ws.Range("L94").Value = 1606.6
$ws.Range("M94").Value = -1492
$ws.Range("N94").Value = -2508.6
$ws.Range("H99").Value = 3251497
$ws.Range("I99").Value = 101546.2
$ws.Range("K99").Value = 101546.2
$ws.Range("M99").Value = -100048.2
$ws.Range("H134").Value = 3275.3809
$ws.Range("J134").Value = 7775.125
$ws.Range("L134").Value = 23325.375
$ws.Range("N134").Value = -28395.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3669.7273
$ws.Range("I31").Value = 2302.111
$ws.Range("K31").Value = 2302.111
$ws.Range("M31").Value = -2007.111
$ws.Range("H34").Value = 3669.7273
$ws.Range("I34").Value = 2302.111
$ws.Range("K34").Value = 2302.111
$ws.Range("M34").Value = -2100.111

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 90
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 90
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 270
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = -550
$ws.Range("H16").Value = 500
$ws.Range("J16").Value = 500
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -1846
$ws.Range("H17").Value = 65
$ws.Range("I17").Value = 70.8
$ws.Range("J17").Value = 57.75
$ws.Range("K17").Value = 212.4
$ws.Range("L17").Value = 173.25
$ws.Range("M17").Value = -43.39999999999998
$ws.Range("N17").Value = -511.25
$ws.Range("H55").Value = 3289.3333
$ws.Range("I55").Value = 934.375
$ws.Range("J55").Value = 7999.25
$ws.Range("K55").Value = 2803.125
$ws.Range("L55").Value = 23997.75
$ws.Range("M55").Value = -2626.125
$ws.Range("N55").Value = -24351.75
$ws.Range("H68").Value = 92857.09
$ws.Range("J68").Value = 101892.8
$ws.Range("L68").Value = 305678.4
$ws.Range("N68").Value = -307300.4
$ws.Range("H71").Value = 92857.09
$ws.Range("J71").Value = 101892.8
$ws.Range("L71").Value = 917035.2000000001
$ws.Range("N71").Value = -925147.2000000001
$ws.Range("H132").Value = 9310.77
$ws.Range("J132").Value = 10028.333
$ws.Range("L132").Value = 90254.997
$ws.Range("N132").Value = -95314.997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8691.462
$ws.Range("I70").Value = 7714
$ws.Range("K70").Value = 7714
$ws.Range("M70").Value = -7444
$ws.Range("H73").Value = 8691.462
$ws.Range("I73").Value = 7714
$ws.Range("K73").Value = 7714
$ws.Range("M73").Value = -6778
$ws.Range("H109").Value = 42244.152
$ws.Range("J109").Value = 45139.5
$ws.Range("L109").Value = 45139.5
$ws.Range("N109").Value = -47219.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2297.818
$ws.Range("I7").Value = 1034.5
$ws.Range("K7").Value = 1034.5
$ws.Range("M7").Value = -922.5
$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 5000
$ws.Range("N29").Value = -5590
$ws.Range("H40").Value = 5054415.5
$ws.Range("I40").Value = 4144.143
$ws.Range("J40").Value = 13892390
$ws.Range("K40").Value = 4144.143
$ws.Range("L40").Value = 13892390
$ws.Range("M40").Value = -4008.143
$ws.Range("N40").Value = -13892662
$ws.Range("H46").Value = 13729.25
$ws.Range("I46").Value = 15490.571
$ws.Range("K46").Value = 15490.571
$ws.Range("M46").Value = -15302.571
$ws.Range("H93").Value = 1056.4615
$ws.Range("I93").Value = 977.8333
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 977.8333
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 270.1667
$ws.Range("N93").Value = -4496
$ws.Range("H100").Value = 7840.2383
$ws.Range("I100").Value = 8263.611000000001
$ws.Range("K100").Value = 8263.611000000001
$ws.Range("M100").Value = -7722.611000000001
$ws.Range("H126").Value = 2297.818
$ws.Range("I126").Value = 1034.5
$ws.Range("K126").Value = 3103.5
$ws.Range("M126").Value = -633.5
$ws.Range("H132").Value = 4819.1113
$ws.Range("I132").Value = 5196
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 15588
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -13058
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 2000
$ws.Range("I34").Value = 2000
$ws.Range("K34").Value = 2000
$ws.Range("M34").Value = -1797
$ws.Range("H75").Value = 30476.191
$ws.Range("J75").Value = 31052.63
$ws.Range("L75").Value = 31052.63
$ws.Range("N75").Value = -32924.63
$ws.Range("H78").Value = 30476.191
$ws.Range("J78").Value = 31052.63
$ws.Range("L78").Value = 93157.89
$ws.Range("N78").Value = -102517.89
$ws.Range("H96").Value = 2400175.2
$ws.Range("J96").Value = 3511738.2
$ws.Range("L96").Value = 3511738.2
$ws.Range("N96").Value = -3514484.2
$ws.Range("H136").Value = 1411.1666
$ws.Range("I136").Value = 1190.7333
$ws.Range("J136").Value = 2513.3333
$ws.Range("K136").Value = 3572.199900000001
$ws.Range("L136").Value = 7539.999899999999
$ws.Range("M136").Value = -1022.199900000001
$ws.Range("N136").Value = -12639.9999
